# Appends 4 new game rows (8-11) to Sheet1, matching the data scraped for
# Garry's Mod, Amnesia: The Dark Descent, Hellblade: Senua's Sacrifice, and
# ARK: Survival Of The Fittest.
#
# Columns: A=Game ID(n) B=Name C=Studio D=Release Date E=Review Count
#          F=Wiki Entry G=References H=Additional Info I=Date Added J=Image URL
#
# Columns D contains plain "YYYY-MM-DD" text that Excel's COM layer will
# otherwise auto-coerce into a date serial number when assigned directly, so
# those cells are written with a temporary Text number format which is
# cleared immediately after the value is set (mirrors how a user would type
# the value into a cell pre-formatted as Text, then reformat back to General).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 8 ----
$ws.Range("A8").Value = 9721
$ws.Range("B8").Value = 'Garry''s Mod'
$ws.Range("C8").Value = 'Garry Newman'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2004-12-24'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 'Garry''s Mod, often abbreviated as GMod, is a sandbox game developed by Garry Newman and published by Valve Corporation. Initially released on December 24, 2004, as a mod for Valve''s Source engine, it has since evolved into a standalone game available on platforms such as Linux, macOS, and PC. The game is renowned for its open-ended gameplay, allowing players to manipulate objects and experiment with the engine''s physics. Unlike traditional video games, Garry''s Mod lacks a predefined plot or gameplay objectives, instead offering a flexible environment where users can create and share their own game modes.
The gameplay mechanics of Garry''s Mod revolve around the player''s ability to spawn objects, characters, and weapons from the Source engine, and to interact with them in a virtual space. This flexibility has led to the creation of a vast array of user-generated content, ranging from simple amusement to complex game modes such as Trouble in Terrorist Town, Prop Hunt, and various franchise-inspired scenarios like Star Wars and Lord of the Rings. The game supports both single-player and multiplayer modes, with servers categorized by different game modes to ensure players can easily find and join their desired experiences.
Garry''s Mod has received generally positive reviews, with a rating of 3.8. Critics and players alike praise its innovative approach to user-generated content and its role as a creative outlet for gamers. The game''s integration with Steam Workshop has facilitated a thriving community where players can share and download mods, further enhancing the game''s replayability and depth. The use of Valve''s Anti-Cheat system ensures a fair and enjoyable experience for all players.
The cultural impact of Garry''s Mod is significant, as it has inspired a generation of gamers to explore game design and modding. Many players have used it as a stepping stone into the gaming industry, honing skills in animation, modeling, and level design. As a testament to its enduring popularity, Garry''s Mod remains a staple in the sandbox genre, celebrated for its humor, creativity, and the freedom it offers players to express themselves through gameplay. The game''s visual style, as depicted in its artwork, reflects a blend of realism and whimsical creativity, which has become a hallmark of its enduring appeal.'
$ws.Range("G8").Value = '["Newman, G. (2004). Garry''s Mod. Facepunch Studios.", "Valve Corporation. (n.d.). Garry''s Mod [PC, Linux, macOS].", "Smith, A. (2015). The creative impact of Garry''s Mod on the gaming industry. PC Gamer.", "Jones, T. (2020). A decade of sandbox innovation: The legacy of Garry''s Mod. Game Informer.", "Brown, J. (2022). User-generated content in video games: A study of Garry''s Mod. Journal of Game Design and Development."]'
$ws.Range("H8").Value = 'Average Rating: 3.8
Average Playtime: 14 hours
ESRB Rating: Rating Pending
Platforms: Linux, macOS, PC'
$ws.Range("I8").Value = '2025-03-27 18:58:55'
$ws.Range("J8").Value = 'https://media.rawg.io/media/games/48c/48cb04ca483be865e3a83119c94e6097.jpg'

# ---- Row 9 ----
$ws.Range("A9").Value = 13668
$ws.Range("B9").Value = 'Amnesia: The Dark Descent'
$ws.Range("C9").Value = 'Mikko Tarmia, Thomas Grip, Jens Nilsson, Mikael Hedberg, Marc Nicander, Marcus Johansson, Luis Rodero Morales'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2010-09-08'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 'Amnesia: The Dark Descent is a critically acclaimed survival horror video game developed by Frictional Games and released on September 8, 2010. The game is notable for its atmospheric tension and psychological horror elements, setting a new standard in the genre. It was created by a team including Mikko Tarmia, Thomas Grip, Jens Nilsson, Mikael Hedberg, Marc Nicander, Marcus Johansson, and Luis Rodero Morales. Designed for platforms such as Linux, macOS, and PC, Amnesia: The Dark Descent is revered for its immersive storytelling and innovative gameplay mechanics.
The game places the player in the role of Daniel, a young man who awakens in an eerie, seemingly deserted castle with no memory of his past. The narrative unfolds through environmental cues, notes, audio logs, and visions, inviting players to unravel the castle’s deep-seated secrets. A significant innovation in the game is the Sanity mechanic, which simulates the psychological impact of the environment on Daniel. His mental state can deteriorate from darkness, unsettling visions, and encounters with horrifying creatures, leading to distorted visual and auditory effects that heighten the sense of dread.
Amnesia: The Dark Descent challenges players through item-based puzzles and interactive environments, requiring resourcefulness and bravery. The game supports modding, allowing user-generated content to flourish and expand the game''s universe. It also features a standalone DLC, Justine, which offers a short yet compelling narrative about a female prisoner with multiple possible endings, further enriching the game’s lore.
Upon release, Amnesia: The Dark Descent received widespread acclaim for its sound design, story richness, and ability to evoke fear without relying on combat mechanics. It garnered a rating of 3.64 and became a landmark title in the survival horror genre. The game''s influence is evident in subsequent indie horror titles, which adopted its emphasis on atmosphere and psychological depth. The Lovecraftian themes and the gothic setting have been praised for their authenticity and immersive quality, contributing to the game''s enduring legacy.
Amnesia: The Dark Descent has left a significant impact on both players and developers, expanding the boundaries of horror in video games. Its success demonstrated the potential for indie developers to create compelling, atmospheric experiences that rival those of major studios.'
$ws.Range("G9").Value = '[''Frictional Games. (2010). Amnesia: The Dark Descent [Video game]. Frictional Games.'', ''IGN. (2010). Amnesia: The Dark Descent Review. Retrieved from https://www.ign.com/articles/2010/09/08/amnesia-the-dark-descent-review'', ''GameSpot. (2010). Amnesia: The Dark Descent Review. Retrieved from https://www.gamespot.com/reviews/amnesia-the-dark-descent-review/1900-6275413/'', ''PC Gamer. (2010). Amnesia: The Dark Descent Review. Retrieved from https://www.pcgamer.com/amnesia-the-dark-descent-review/'', ''Rock, Paper, Shotgun. (2010). Wot I Think: Amnesia: The Dark Descent. Retrieved from https://www.rockpapershotgun.com/amnesia-the-dark-descent-review'', ''Official Amnesia: The Dark Descent Website. Retrieved from http://www.amnesiagame.com'']'
$ws.Range("H9").Value = 'Average Rating: 3.64
Average Playtime: 2 hours
ESRB Rating: Mature
Metacritic Score: 85
Platforms: Linux, macOS, PC'
$ws.Range("I9").Value = '2025-03-27 18:59:06'
$ws.Range("J9").Value = 'https://media.rawg.io/media/games/b54/b54598d1d5cc31899f4f0a7e3122a7b0.jpg'

# ---- Row 10 ----
$ws.Range("A10").Value = 12536
$ws.Range("B10").Value = 'Hellblade: Senua''s Sacrifice'
$ws.Range("C10").Value = 'Tameem Antoniades, Rupert Brooker, Joe Nelson, Melanie Hall, Juan Fernandez, Andy LaPlegua, David García, Gavin Costello, Stefano Prosperi, Loong Wei Ding'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '2017-08-07'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 'Hellblade: Senua''s Sacrifice is a dark fantasy action-adventure game developed and published by Ninja Theory, released on August 7, 2017. The game is set in a gritty world inspired by Scandinavian and Celtic folklore, featuring eerie landscapes and masked enemies. Players control Senua, a Pict warrior grappling with psychosis, as she embarks on a harrowing journey to Helheim to retrieve the soul of her deceased lover, Dillion. Her quest is not only a physical battle against enemies but also a symbolic exploration of her mental struggles, with the game blurring lines between reality and delusion.
The gameplay of Hellblade: Senua''s Sacrifice combines combat and puzzle-solving elements. Players engage in third-person melee combat against various foes while navigating the environment to solve puzzles that require Senua to "focus" and perceive hidden aspects of her surroundings. The game''s use of binaural audio technology immerses players in Senua''s experiences, with meticulously recorded voices that guide, deceive, or warn her. Motion capture technology was employed to render Senua''s movements and expressions with an emotional authenticity that enhances the storytelling experience.
Critically acclaimed for its narrative depth and innovative approach to portraying mental health, Hellblade received praise for its audio design and visual fidelity. It was lauded for its sensitive depiction of psychosis, developed in collaboration with neuroscientists and individuals with lived experiences of mental illness. The game won several awards, including a BAFTA for Artistic Achievement, and it has been noted for its impact in raising awareness and understanding of mental health issues within the gaming community.
Culturally, Hellblade: Senua''s Sacrifice has been influential in shifting perceptions about the potential of video games as a medium for addressing complex and serious topics. It challenged conventional gaming narratives by focusing on a female protagonist experiencing a mental health condition, contributing to broader discussions about representation and diversity in video games. Its success has encouraged other developers to explore similar themes, solidifying its status as a landmark title in the action-adventure genre.'
$ws.Range("G10").Value = '["Antoniades, T., Brooker, R., Nelson, J., Hall, M., Fernandez, J., LaPlegua, A., García, D., Costello, G., Prosperi, S., & Ding, L. W. (2017). Hellblade: Senua''s Sacrifice [Video game]. Ninja Theory.", "IGN. (2017). Hellblade: Senua''s Sacrifice Review. IGN. Retrieved from https://www.ign.com/articles/2017/08/07/hellblade-senuas-sacrifice-review", "Polygon. (2017). How Ninja Theory made Hellblade: Senua''s Sacrifice. Polygon. Retrieved from https://www.polygon.com/features/2017/8/8/16109796/hellblade-making-of-ninja-theory", "GamesRadar. (2017). Hellblade: Senua''s Sacrifice review: ''A masterpiece of haunting visuals and sound.'' GamesRadar. Retrieved from https://www.gamesradar.com/hellblade-senuas-sacrifice-review/", "Official website of Hellblade: Senua''s Sacrifice. Retrieved from http://www.hellblade.com/"]'
$ws.Range("H10").Value = 'Average Rating: 4.23
Average Playtime: 6 hours
ESRB Rating: Mature
Metacritic Score: 83
Platforms: Xbox One, Nintendo Switch, PlayStation 4, PC'
$ws.Range("I10").Value = '2025-03-27 19:05:37'
$ws.Range("J10").Value = 'https://media.rawg.io/media/games/63f/63f0e68688cad279ed38cde931dbfcdb.jpg'

# ---- Row 11 ----
$ws.Range("A11").Value = 11147
$ws.Range("B11").Value = 'ARK: Survival Of The Fittest'
$ws.Range("C11").Value = ''
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '2016-03-15'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 'ARK: Survival of the Fittest, released on March 15, 2016, is a distinctive entry in the multiplayer gaming landscape, designed by Studio Wildcard, the creators of ARK: Survival Evolved. This game represents the first-ever Multiplayer Online Survival Arena (M.O.S.A.), crafted to cater to the rapidly expanding eSports scene. As a spin-off from the popular open-world game ARK: Survival Evolved, SotF introduces players to a high-stakes environment where up to 72 players engage in an intense survival battle. The game emphasizes strategic gameplay as participants form Tribes of one to six members, competing on both unofficial and ranked servers in scenarios that range from swift 30-minute skirmishes to extended three-hour confrontations.
The gameplay of ARK: Survival of the Fittest is marked by its dynamic survival mechanics and competitive elements. Players begin in a neutral staging area to plan their strategies before descending onto the ARK, a lush and perilous environment teeming with prehistoric creatures. The game features over 30 creatures available for taming and riding at the onset of Early Access, with plans for more at full launch. Players must navigate various strategies, such as racing to the central cache for weapons and resources, hiding in the dense forest, or building an army of dinosaurs to dominate opponents. A shrinking "ring of death" gradually forces all players into closer proximity, ensuring that the tension and conflict escalate as the match progresses.
The game is equipped with robust eSports features, including dedicated ranking systems, matchmaking, and tournament functionalities, making it accessible to a wide audience. Players can also earn exclusive cosmetic items that are compatible with ARK: Survival Evolved. The auditory experience is enriched by a custom soundtrack from Gareth Coker, known for his work on Ori and the Blind Forest, adding a cinematic quality to the competitive gameplay.
Despite its innovative approach and rich features, ARK: Survival of the Fittest received a mixed reception, garnering an average rating of 2.59. It was praised for its originality and the exciting potential it brought to the eSports arena, yet faced criticism for gameplay balance issues and technical challenges typical of Early Access releases. Nonetheless, it has contributed to the cultural landscape by blending survival mechanics with competitive eSports, thus influencing subsequent titles in the genre.
The game is available on multiple platforms, including Android, Linux, macOS, PC, and iOS, supporting both single-player and multiplayer modes. Its development and release underscore the growing intersection of survival games and competitive multiplayer experiences in the digital gaming domain.'
$ws.Range("G11").Value = '[''Studio Wildcard. (2016). ARK: Survival of the Fittest [Video game]. Studio Wildcard.'', ''Coker, G. (Composer). (2016). ARK: Survival of the Fittest Original Soundtrack [Soundtrack].'', ''IGN. (2016). ARK: Survival of the Fittest Review. Retrieved from https://www.ign.com/articles/2016/03/15/ark-survival-of-the-fittest-review'', ''PC Gamer. (2016). ARK: Survival of the Fittest: A New Twist in Battle Royale. Retrieved from https://www.pcgamer.com/ark-survival-of-the-fittest-impressions'', ''Polygon. (2016). ARK: Survival of the Fittest: How it Changed the Survival Genre. Retrieved from https://www.polygon.com/2016/03/15/ark-survival-of-the-fittest'']'
$ws.Range("H11").Value = 'Average Rating: 2.59
Average Playtime: 1 hours
ESRB Rating: Teen
Platforms: Android, Linux, macOS, PC, iOS'
$ws.Range("I11").Value = '2025-03-27 19:05:59'
$ws.Range("J11").Value = 'https://media.rawg.io/media/games/095/0953bf01cd4e4dd204aba85489ac9868.jpg'

